$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.739.09'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').Value = '1.626.46'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '214.82'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').Value = '0.5099'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.2559'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.06316'
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('D10').Value = '19.37'
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.631.60'
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.224'
$ws.Range('E13').Value = '  -1.30%  '
$ws.Range('D14').Value = '1.850.13'
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('D15').Value = '0.5506'
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('D16').Value = '63.49'
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('D17').Value = '0.0₅7486'
$ws.Range('E17').Value = '  -2.53%  '
$ws.Range('D18').Value = '25.767.43'
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '4.400'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '193.60'
$ws.Range('E21').Value = '  -2.57%  '
$ws.Range('D22').Value = '9.796'
$ws.Range('E22').Value = '  -1.02%  '
$ws.Range('D23').Value = '5.986'
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('D24').Value = '1.003'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('D26').Value = '141.91'
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('D27').Value = '0.1254'
$ws.Range('E27').Value = '  +5.43%  '
$ws.Range('D28').Value = '15.49'
$ws.Range('E28').Value = '  -0.90%  '
$ws.Range('D29').Value = '6.722'
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('D30').Value = '1.236'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = '0.04873'
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').Value = '3.233'
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('D33').Value = '3.150'
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('D34').Value = '1.535'
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('D35').Value = '2.376'
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('D36').Value = '0.8908'
$ws.Range('E36').Value = '  -1.24%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').Value = '2.539'
$ws.Range('E37').Value = '  -1.77%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '0.5494'
$ws.Range('E38').Value = '  +0.95%  '
$ws.Range('D39').Value = '1.110.63'
$ws.Range('E39').Value = '  -2.91%  '
$ws.Range('D40').Value = '0.01544'
$ws.Range('E40').Value = '  -1.10%  '
$ws.Range('D41').Value = '1.000'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('D42').Value = '5.546'
$ws.Range('E42').Value = '  +2.62%  '
$ws.Range('D43').Value = '0.7965'
$ws.Range('E43').Value = '  -1.77%  '
$ws.Range('D44').Value = '97.14'
$ws.Range('E44').Value = '  -2.18%  '
$ws.Range('D45').Value = '1.773.80'
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('E46').Value = '  -11.19%  '
$ws.Range('D47').Value = '0.4432'
$ws.Range('E47').Value = '  -2.20%  '
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('D49').Value = '54.52'
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('D50').Value = '0.05135'
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').Value = '7.502'
$ws.Range('E51').Value = '  +2.79%  '
